# [add] poi import excel template
#
# Adds a "签名" / "${signName}" label-value pair on row 5 (columns F:G),
# styled the same as the existing F3/G3 header cells, and moves the
# worksheet's active selection from C1:D1 down to the single cell C10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 5: F5 = "签名", G5 = "${signName}"
$ws.Range("F5").Value = "签名"
$ws.Range("G5").Value = '${signName}'

# Match the formatting already used by the F3/G3 header cells (centered,
# size-18 font) so the new row reuses the existing cell style.
$ws.Range("F5:G5").HorizontalAlignment = $ws.Range("F3:G3").HorizontalAlignment
$ws.Range("F5:G5").VerticalAlignment = $ws.Range("F3:G3").VerticalAlignment
$ws.Range("F5:G5").Font.Size = $ws.Range("F3:G3").Font.Size

# Move the active selection to the single cell C10.
$ws.Range("C10").Select()
